$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "StatQuery" column (C) on rows 2-4 held a single-line Cypher query.
# It is being replaced everywhere by a reformatted, multi-line version of
# the same query (line breaks + the OPTIONAL MATCH clauses split in two).
$newStatQuery = @'
MATCH (s:study)
  WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies
  MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies
  MATCH (d:diagnosis)
  WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies
  MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
    WHERE demo.breed IN ['Labrador Retriever']
  OPTIONAL MATCH (f:file)-[*]->(c)
  OPTIONAL MATCH (samp:sample)-[*]->(c)
  WITH DISTINCT c AS c, p, s, demo, diag, f, samp
  RETURN count(DISTINCT(f)) as number_of_files ,
             count(DISTINCT(samp)) as number_of_sample ,
             count(DISTINCT(c.case_id)) as number_of_cases ,
             count(DISTINCT(s.clinical_study_designation)) as number_of_study
'@

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Row heights re-wrap to the new text (wrapText is on for these cells).
# Row 2 and row 4 keep their previous auto-computed heights; row 3 grows.
$ws.Rows.Item(2).RowHeight = 275.5
$ws.Rows.Item(3).RowHeight = 246.5
$ws.Rows.Item(4).RowHeight = 246.5

# Selection/view moved from C2 to B2, scrolled down one row.
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B2").Select() | Out-Null
